$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/week number, report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/7/2025  Through  4/13/2025"

# --- Crime Complaints table updates ---
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = 0
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J15").Value = 6
$ws.Range("J15").NumberFormat = "#,##0"
$ws.Range("K15").Value = 50
$ws.Range("K15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L15").Value = 50
$ws.Range("L15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C16").Value = 1
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 8
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").Value = -87.5
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F16").Value = 11
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("G16").Value = 20
$ws.Range("G16").NumberFormat = "#,##0"
$ws.Range("H16").Value = -45
$ws.Range("H16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I16").Value = 32
$ws.Range("I16").NumberFormat = "#,##0"
$ws.Range("J16").Value = 61
$ws.Range("J16").NumberFormat = "#,##0"
$ws.Range("K16").Value = -47.540983606557
$ws.Range("K16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L16").Value = -41.818181818181
$ws.Range("L16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M16").Value = -50
$ws.Range("M16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N16").Value = -88.811188811188
$ws.Range("N16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C17").Value = 4
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("E17").Value = 0
$ws.Range("E17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F17").Value = 19
$ws.Range("F17").NumberFormat = "#,##0"
$ws.Range("G17").Value = 17
$ws.Range("G17").NumberFormat = "#,##0"
$ws.Range("H17").Value = 11.764705882352
$ws.Range("H17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I17").Value = 76
$ws.Range("I17").NumberFormat = "#,##0"
$ws.Range("J17").Value = 73
$ws.Range("J17").NumberFormat = "#,##0"
$ws.Range("K17").Value = 4.109589041095
$ws.Range("K17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L17").Value = 7.042253521126
$ws.Range("L17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M17").Value = 26.666666666666
$ws.Range("M17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N17").Value = 11.764705882352
$ws.Range("N17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C18").Value = 9
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 3
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("E18").Value = 200
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F18").Value = 24
$ws.Range("F18").NumberFormat = "#,##0"
$ws.Range("G18").Value = 16
$ws.Range("G18").NumberFormat = "#,##0"
$ws.Range("H18").Value = 50
$ws.Range("H18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I18").Value = 58
$ws.Range("I18").NumberFormat = "#,##0"
$ws.Range("J18").Value = 61
$ws.Range("J18").NumberFormat = "#,##0"
$ws.Range("K18").Value = -4.918032786885
$ws.Range("K18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L18").Value = 9.43396226415
$ws.Range("L18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M18").Value = -51.666666666666
$ws.Range("M18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N18").Value = -89.624329159212
$ws.Range("N18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C19").Value = 15
$ws.Range("C19").NumberFormat = "#,##0"
$ws.Range("D19").Value = 12
$ws.Range("D19").NumberFormat = "#,##0"
$ws.Range("E19").Value = 25
$ws.Range("E19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F19").Value = 57
$ws.Range("F19").NumberFormat = "#,##0"
$ws.Range("G19").Value = 53
$ws.Range("G19").NumberFormat = "#,##0"
$ws.Range("H19").Value = 7.54716981132
$ws.Range("H19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I19").Value = 158
$ws.Range("I19").NumberFormat = "#,##0"
$ws.Range("J19").Value = 212
$ws.Range("J19").NumberFormat = "#,##0"
$ws.Range("K19").Value = -25.471698113207
$ws.Range("K19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L19").Value = -8.670520231213
$ws.Range("L19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M19").Value = 49.056603773584
$ws.Range("M19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N19").Value = 0
$ws.Range("N19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C20").Value = 9
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 6
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("E20").Value = 50
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F20").Value = 27
$ws.Range("F20").NumberFormat = "#,##0"
$ws.Range("G20").Value = 27
$ws.Range("G20").NumberFormat = "#,##0"
$ws.Range("H20").Value = 0
$ws.Range("H20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I20").Value = 69
$ws.Range("I20").NumberFormat = "#,##0"
$ws.Range("J20").Value = 98
$ws.Range("J20").NumberFormat = "#,##0"
$ws.Range("K20").Value = -29.591836734693
$ws.Range("K20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L20").Value = -20.689655172413
$ws.Range("L20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M20").Value = -26.595744680851
$ws.Range("M20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N20").Value = -93.441064638783
$ws.Range("N20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C21").Value = 38
$ws.Range("C21").NumberFormat = "#,##0"
$ws.Range("D21").Value = 34
$ws.Range("D21").NumberFormat = "#,##0"
$ws.Range("E21").Value = 11.764705882352
$ws.Range("E21").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F21").Value = 139
$ws.Range("F21").NumberFormat = "#,##0"
$ws.Range("G21").Value = 135
$ws.Range("G21").NumberFormat = "#,##0"
$ws.Range("H21").Value = 2.962962962962
$ws.Range("H21").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I21").Value = 402
$ws.Range("I21").NumberFormat = "#,##0"
$ws.Range("J21").Value = 512
$ws.Range("J21").NumberFormat = "#,##0"
$ws.Range("K21").Value = -21.484375
$ws.Range("K21").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L21").Value = -10.067114093959
$ws.Range("L21").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M21").Value = -10.467706013363
$ws.Range("M21").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N21").Value = -81.179775280898
$ws.Range("N21").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 2
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("G22").Value = 2
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("H22").Value = 0
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I22").Value = 5
$ws.Range("I22").NumberFormat = "#,##0"
$ws.Range("K22").Value = -37.5
$ws.Range("K22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L22").Value = -28.571428571428
$ws.Range("L22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M22").Value = -16.666666666666
$ws.Range("M22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C24").Value = 30
$ws.Range("C24").NumberFormat = "#,##0"
$ws.Range("D24").Value = 26
$ws.Range("D24").NumberFormat = "#,##0"
$ws.Range("E24").Value = 15.384615384615
$ws.Range("E24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F24").Value = 112
$ws.Range("F24").NumberFormat = "#,##0"
$ws.Range("G24").Value = 94
$ws.Range("G24").NumberFormat = "#,##0"
$ws.Range("H24").Value = 19.148936170212
$ws.Range("H24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I24").Value = 345
$ws.Range("I24").NumberFormat = "#,##0"
$ws.Range("J24").Value = 406
$ws.Range("J24").NumberFormat = "#,##0"
$ws.Range("K24").Value = -15.024630541871
$ws.Range("K24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L24").Value = -7.008086253369
$ws.Range("L24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M24").Value = 13.486842105263
$ws.Range("M24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C25").Value = 12
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("D25").Value = 17
$ws.Range("D25").NumberFormat = "#,##0"
$ws.Range("E25").Value = -29.411764705882
$ws.Range("E25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F25").Value = 39
$ws.Range("F25").NumberFormat = "#,##0"
$ws.Range("G25").Value = 40
$ws.Range("G25").NumberFormat = "#,##0"
$ws.Range("H25").Value = -2.5
$ws.Range("H25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I25").Value = 140
$ws.Range("I25").NumberFormat = "#,##0"
$ws.Range("J25").Value = 131
$ws.Range("J25").NumberFormat = "#,##0"
$ws.Range("K25").Value = 6.870229007633
$ws.Range("K25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L25").Value = 15.702479338843
$ws.Range("L25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C26").Value = 7
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 9
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = -22.222222222222
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 42
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("G26").Value = 49
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("H26").Value = -14.285714285714
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I26").Value = 130
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("J26").Value = 161
$ws.Range("J26").NumberFormat = "#,##0"
$ws.Range("K26").Value = -19.254658385093
$ws.Range("K26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L26").Value = -6.474820143884
$ws.Range("L26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M26").Value = -37.5
$ws.Range("M26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 1
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("H27").Value = 0
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J27").Value = 7
$ws.Range("J27").NumberFormat = "#,##0"
$ws.Range("K27").Value = 42.857142857142
$ws.Range("K27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L27").Value = 25
$ws.Range("L27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D28").Value = 3
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -66.666666666666
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F28").Value = 5
$ws.Range("F28").NumberFormat = "#,##0"
$ws.Range("G28").Value = 5
$ws.Range("G28").NumberFormat = "#,##0"
$ws.Range("H28").Value = 0
$ws.Range("H28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I28").Value = 15
$ws.Range("I28").NumberFormat = "#,##0"
$ws.Range("J28").Value = 11
$ws.Range("J28").NumberFormat = "#,##0"
$ws.Range("K28").Value = 36.363636363636
$ws.Range("K28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L28").Value = -21.052631578947
$ws.Range("L28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "0"
$ws.Range("H31").Value = -100
$ws.Range("H31").NumberFormat = "#,##0.0;""-""#,##0.0"
